# Add a new "time_taken" column (F) to the panel worksheet, recording the
# time each row's data was retrieved from panelApp, for metadata purposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1 -------------------------------------------------------
# Copy the formatting from the neighbouring header cell (E1, which is bold,
# bordered & centered) onto F1, then set its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# --- Data cells F2:F46 -----------------------------------------------------
$timeTaken = @(
    "2021-10-05 10:52:05.095087",
    "2021-10-05 10:52:05.095099",
    "2021-10-05 10:52:05.095102",
    "2021-10-05 10:52:05.095105",
    "2021-10-05 10:52:05.095108",
    "2021-10-05 10:52:05.095110",
    "2021-10-05 10:52:05.095113",
    "2021-10-05 10:52:05.095116",
    "2021-10-05 10:52:05.095119",
    "2021-10-05 10:52:05.095122",
    "2021-10-05 10:52:05.095125",
    "2021-10-05 10:52:05.095127",
    "2021-10-05 10:52:05.095130",
    "2021-10-05 10:52:05.095132",
    "2021-10-05 10:52:05.095135",
    "2021-10-05 10:52:05.095138",
    "2021-10-05 10:52:05.095140",
    "2021-10-05 10:52:05.095143",
    "2021-10-05 10:52:05.095146",
    "2021-10-05 10:52:05.095149",
    "2021-10-05 10:52:05.095151",
    "2021-10-05 10:52:05.095154",
    "2021-10-05 10:52:05.095156",
    "2021-10-05 10:52:05.095159",
    "2021-10-05 10:52:05.095162",
    "2021-10-05 10:52:05.095165",
    "2021-10-05 10:52:05.095168",
    "2021-10-05 10:52:05.095170",
    "2021-10-05 10:52:05.095173",
    "2021-10-05 10:52:05.095175",
    "2021-10-05 10:52:05.095178",
    "2021-10-05 10:52:05.095181",
    "2021-10-05 10:52:05.095184",
    "2021-10-05 10:52:05.095187",
    "2021-10-05 10:52:05.095189",
    "2021-10-05 10:52:05.095192",
    "2021-10-05 10:52:05.095194",
    "2021-10-05 10:52:05.095197",
    "2021-10-05 10:52:05.095199",
    "2021-10-05 10:52:05.095202",
    "2021-10-05 10:52:05.095205",
    "2021-10-05 10:52:05.095208",
    "2021-10-05 10:52:05.095210",
    "2021-10-05 10:52:05.095213",
    "2021-10-05 10:52:05.095215"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "Added time_taken column (F1:F46)."
